# Update "想去人数" (interested count) figures on the 展览 and 全部类型 sheets
# to reflect newly generated data (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览"   = @{ F6 = 3805; F16 = 232; F19 = 96; F21 = 3342; F22 = 5680; F43 = 28 }
    "全部类型" = @{ F6 = 3805; F17 = 232; F20 = 96; F22 = 3342; F23 = 5680; F44 = 28 }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
